# Generate Report for Handback
# Update the handoff/handback datetime stamps for the first file row
# (21aa40b0-7c4a-4773-a4dd-feb631c62cdb...) on both the zh-cn and de-de
# language sheets, reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-17 09:07:52"
$wsZhCn.Range("G2").Value = "2016-02-17 09:08:37"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-17 09:08:05"
$wsDeDe.Range("G2").Value = "2016-02-17 09:08:55"
